# Insert a new data row for "Fruta / hortaliza, semanal" update.
# A new record is inserted at row 37, pushing the existing rows 37..156
# down to 38..157 (matching the sheet's weekly/periodic logic of
# prepending the newest observation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Physically insert a new row at position 37; this shifts all rows
# 37..156 down to 38..157 (carrying their existing formatting/styles).
$ws.Rows(37).Insert()

# Populate the newly inserted row 37 with the new record.
$ws.Cells.Item(37, 1).Value = 7
$ws.Cells.Item(37, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(37, 3).Value = "Ñuble"
$ws.Cells.Item(37, 4).Value = 45114
$ws.Cells.Item(37, 5).Value = 16
$ws.Cells.Item(37, 6).Value = "Fruta"
$ws.Cells.Item(37, 7).Value = 100108
$ws.Cells.Item(37, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(37, 9).Value = 100108002
$ws.Cells.Item(37, 10).Value = "Mango"
$ws.Cells.Item(37, 11).Value = "Sin especificar"
$ws.Cells.Item(37, 12).Value = "Primera"
$ws.Cells.Item(37, 13).Value = 40
$ws.Cells.Item(37, 14).Value = 9000
$ws.Cells.Item(37, 15).Value = 10000
$ws.Cells.Item(37, 16).Value = 9500
$ws.Cells.Item(37, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(37, 18).Value = "Brasil"
$ws.Cells.Item(37, 19).Value = 2375
$ws.Cells.Item(37, 20).Value = 4
